# Update gh-pages to output generated at 456a3b4
# Applies the "杭州-漫展信息.xlsx" data refresh:
#   - refresh "想去人数" (F column) counters on sheet "展览"
#   - insert a newly-scraped event ("杭州·AD05动漫展") as row 30 on sheet "展览",
#     pushing the previous last row ("华盟次元嘉年华...") down to row 31
#   - refresh the same "想去人数" (F column) counters on sheet "全部类型"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------
# 1. "展览" sheet — refresh F-column (想去人数) counts for existing rows
# ---------------------------------------------------------------------
$sheet1Updates = @{
    2  = 241
    3  = 1406
    6  = 30
    7  = 1241
    8  = 1576
    11 = 2273
    12 = 455
    13 = 121
    16 = 90
    17 = 82
    18 = 6200
    19 = 48
    20 = 6076
    21 = 10058
    23 = 174
    24 = 183
    25 = 276
    26 = 499
    27 = 169
    28 = 149
    29 = 4393
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------
# 2. "展览" sheet — insert the new event as row 30.
#    Copy the current last data row (row 30) down to row 31 first (this
#    preserves that row's styling), then overwrite row 30 with the new
#    event's data.
# ---------------------------------------------------------------------
$ws1.Rows.Item(30).Copy() | Out-Null
$ws1.Rows.Item(30).Insert() | Out-Null

# The Insert() above doesn't carry the border formatting of column A over
# to the freshly-inserted row — restore it so A30 matches the sheet's
# existing "index" column styling.
$ws1.Cells.Item(30, 1).Borders.LineStyle = 1

# Row 31 now holds a duplicate of the old row 30 ("华盟次元嘉年华...").
# Fix its sequence number and its "想去人数" count.
$ws1.Cells.Item(31, 1).Value = 30
$ws1.Cells.Item(31, 6).Value = 386

# Row 30 becomes the newly scraped event. Column A (index number) keeps
# the value 29 that the copy already placed there.
#
# B30 looks like a date ("2024-12-14"), so a plain .Value assignment would
# get auto-coerced into a real date serial by Excel's text-to-date
# detection — but the source data stores it as plain text (matching every
# other row in this column). Force text formatting, assign, then restore
# default (General/no-special-format) styling by pulling the format from
# an existing plain-text date cell in the same column.
$ws1.Cells.Item(30, 2).NumberFormat = "@"
$ws1.Cells.Item(30, 2).Value = "2024-12-14"
$ws1.Cells.Item(29, 2).Copy() | Out-Null
$ws1.Cells.Item(30, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Cells.Item(30, 3).Value = "杭州·AD05动漫展"
$ws1.Cells.Item(30, 4).Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws1.Cells.Item(30, 5).Value = "2024.12.14 10:00-12.15 17:00"
$ws1.Cells.Item(30, 6).Value = 42
$ws1.Cells.Item(30, 7).Value = 75
$ws1.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93152"
$ws1.Cells.Item(30, 9).Value = "//i2.hdslb.com/bfs/openplatform/202409/eXuHZu841725265314495.jpeg"

# ---------------------------------------------------------------------
# 3. "全部类型" sheet — refresh F-column (想去人数) counts
# ---------------------------------------------------------------------
$sheet4Updates = @{
    4  = 241
    5  = 1406
    9  = 30
    10 = 1241
    12 = 1576
    15 = 2274
    17 = 455
    18 = 121
    22 = 90
    23 = 82
    24 = 6200
    25 = 48
    26 = 6076
    27 = 10058
    30 = 174
    31 = 183
    32 = 276
    34 = 499
    38 = 169
    39 = 149
    40 = 4393
    46 = 386
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
